$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 374.5
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 374.5
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 374.5
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -714.5

$ws.Range("H19").Value = 1749.6
$ws.Range("I19").Value = 2375
$ws.Range("J19").Value = 1332.6666
$ws.Range("K19").Value = 2375
$ws.Range("L19").Value = 1332.6666
$ws.Range("M19").Value = -2200
$ws.Range("N19").Value = -1682.6666

$ws.Range("H92").Value = 648.8
$ws.Range("I92").Value = 561
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 561
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = 687
$ws.Range("N92").Value = -3496

$ws.Range("H98").Value = 3086.6667
$ws.Range("I98").Value = 2304.7
$ws.Range("J98").Value = 6996.5
$ws.Range("K98").Value = 2304.7
$ws.Range("L98").Value = 6996.5
$ws.Range("M98").Value = -806.6999999999998
$ws.Range("N98").Value = -9992.5

$ws.Range("H122").Value = 3086.6667
$ws.Range("I122").Value = 2304.7
$ws.Range("J122").Value = 6996.5
$ws.Range("K122").Value = 6914.099999999999
$ws.Range("L122").Value = 20989.5
$ws.Range("M122").Value = -4464.099999999999
$ws.Range("N122").Value = -25889.5

$ws.Range("H125").Value = 533
$ws.Range("I125").Value = 533
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 4797
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -2337
$ws.Range("N125").ClearContents()

$ws.Range("H131").Value = 3729
$ws.Range("I131").Value = 3253.8
$ws.Range("J131").Value = 6105
$ws.Range("K131").Value = 9761.400000000001
$ws.Range("L131").Value = 18315
$ws.Range("M131").Value = -4721.400000000001
$ws.Range("N131").Value = -28395

$ws.Range("H141").Value = 811.55554
$ws.Range("I141").Value = 727.7143
$ws.Range("K141").Value = 2183.1429
$ws.Range("M141").Value = 2996.8571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5390.2
$ws.Range("I32").Value = 5859.778
$ws.Range("K32").Value = 5859.778
$ws.Range("M32").Value = -5572.778

$ws.Range("H74").Value = 1479.3334
$ws.Range("I74").Value = 1175.2
$ws.Range("K74").Value = 1175.2
$ws.Range("M74").Value = -301.2

$ws.Range("H77").Value = 1479.3334
$ws.Range("I77").Value = 1175.2
$ws.Range("K77").Value = 5876
$ws.Range("M77").Value = -1508

$ws.Range("H88").Value = 539
$ws.Range("I88").Value = 562.5
$ws.Range("J88").Value = 520.2
$ws.Range("K88").Value = 562.5
$ws.Range("L88").Value = 520.2
$ws.Range("M88").Value = -156.5
$ws.Range("N88").Value = -1332.2

$ws.Range("H91").Value = 539
$ws.Range("I91").Value = 562.5
$ws.Range("J91").Value = 520.2
$ws.Range("K91").Value = 562.5
$ws.Range("L91").Value = 520.2
$ws.Range("M91").Value = 841.5
$ws.Range("N91").Value = -3328.2

$ws.Range("H96").Value = 29961.125
$ws.Range("J96").Value = 29961.125
$ws.Range("L96").Value = 29961.125
$ws.Range("N96").Value = -35453.125

$ws.Range("H122").Value = 1833
$ws.Range("I122").Value = 1899.6
$ws.Range("K122").Value = 5698.799999999999
$ws.Range("M122").Value = -3248.799999999999

$ws.Range("H132").Value = 2401.6667
$ws.Range("I132").Value = 2401.6667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7205.000100000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4675.000100000001
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 638.5
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -150

$ws.Range("H31").Value = 3488.4
$ws.Range("I31").Value = 3000
$ws.Range("K31").Value = 3000
$ws.Range("M31").Value = -2705

$ws.Range("H34").Value = 3488.4
$ws.Range("I34").Value = 3000
$ws.Range("K34").Value = 3000
$ws.Range("M34").Value = -2798

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").ClearContents()

$ws.Range("H119").Value = 43999
$ws.Range("J119").Value = 43999
$ws.Range("L119").Value = 43999
$ws.Range("N119").Value = -53675

$ws.Range("H124").Value = 111748.25
$ws.Range("J124").Value = 116331.336
$ws.Range("L124").Value = 116331.336
$ws.Range("N124").Value = -121241.336

$ws.Range("H132").Value = 3311.1667
$ws.Range("I132").Value = 2466.75
$ws.Range("K132").Value = 7400.25
$ws.Range("M132").Value = -4870.25

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 513.61536
$ws.Range("I113").Value = 528.2857
$ws.Range("J113").Value = 496.5
$ws.Range("K113").Value = 1584.8571
$ws.Range("L113").Value = 1489.5
$ws.Range("M113").Value = 585.1428999999998
$ws.Range("N113").Value = -5829.5

$ws.Range("H138").Value = 1500
$ws.Range("I138").Value = 1500
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 4500
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 640
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 1252299.8
$ws.Range("I22").Value = 1668066.4
$ws.Range("K22").Value = 1668066.4
$ws.Range("M22").Value = -1667537.4

$ws.Range("H102").Value = 2574.6
$ws.Range("I102").Value = 2597
$ws.Range("J102").Value = 2485
$ws.Range("K102").Value = 2597
$ws.Range("L102").Value = 2485
$ws.Range("M102").Value = -975
$ws.Range("N102").Value = -5729

$ws.Range("H122").Value = 855.8
$ws.Range("I122").Value = 783.3333
$ws.Range("K122").Value = 2349.9999
$ws.Range("M122").Value = 100.0001000000002

$ws.Range("H132").Value = 3738.75
$ws.Range("I132").Value = 3651.8333
$ws.Range("K132").Value = 10955.4999
$ws.Range("M132").Value = -8425.499899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 5947
$ws.Range("I9").Value = 1420.5
$ws.Range("J9").Value = 15000
$ws.Range("K9").Value = 1420.5
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = -1196.5
$ws.Range("N9").Value = -15448

$ws.Range("H13").Value = 4997.5
$ws.Range("J13").Value = 4997.5
$ws.Range("L13").Value = 4997.5
$ws.Range("N13").Value = -5277.5

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

$ws.Range("H22").Value = 531.5
$ws.Range("I22").Value = 489.8
$ws.Range("J22").Value = 740
$ws.Range("K22").Value = 489.8
$ws.Range("L22").Value = 740
$ws.Range("M22").Value = -194.8
$ws.Range("N22").Value = -1330

$ws.Range("H27").Value = 531.5
$ws.Range("I27").Value = 489.8
$ws.Range("J27").Value = 740
$ws.Range("K27").Value = 489.8
$ws.Range("L27").Value = 740
$ws.Range("M27").Value = -382.8
$ws.Range("N27").Value = -954

$ws.Range("H30").Value = 694.7143
$ws.Range("I30").Value = 465.75
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 465.75
$ws.Range("L30").Value = 1000
$ws.Range("M30").Value = -357.75
$ws.Range("N30").Value = -1216

$ws.Range("H40").Value = 3499.75
$ws.Range("I40").Value = 3333
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 3333
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -3197
$ws.Range("N40").Value = -4272

$ws.Range("H61").Value = 5752
$ws.Range("I61").Value = 5003
$ws.Range("K61").Value = 5003
$ws.Range("M61").Value = -4801

$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H113").Value = 5752
$ws.Range("I113").Value = 5003
$ws.Range("K113").Value = 5003
$ws.Range("M113").Value = -2833

$ws.Range("H122").Value = 5760.4
$ws.Range("I122").Value = 5760.4
$ws.Range("K122").Value = 17281.2
$ws.Range("M122").Value = -14831.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1887

$ws.Range("H9").Value = 17502.5
$ws.Range("I9").Value = 17502.5
$ws.Range("K9").Value = 17502.5
$ws.Range("M9").Value = -17362.5

$ws.Range("H18").Value = 1500
$ws.Range("J18").Value = 1500
$ws.Range("L18").Value = 1500
$ws.Range("N18").Value = -1846

$ws.Range("H20").Value = 100011
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 100011
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 100011
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -100491

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H81").Value = 1679.8
$ws.Range("I81").Value = 483
$ws.Range("K81").Value = 966
$ws.Range("M81").Value = 95

$ws.Range("H84").Value = 1679.8
$ws.Range("I84").Value = 483
$ws.Range("K84").Value = 4830
$ws.Range("M84").Value = 474

$ws.Range("H98").Value = 19000
$ws.Range("J98").Value = 19000
$ws.Range("L98").Value = 19000
$ws.Range("N98").Value = -24990

$ws.Range("H122").Value = 386
$ws.Range("I122").Value = 386
$ws.Range("K122").Value = 1158
$ws.Range("M122").Value = 1292
